# 815 - Phase accordion
# The "Mandatory" flag (column I) on the Events sheet defaulted to 1 (true)
# for every event row; the template is updated so events are no longer
# mandatory by default (value 0), except for row 5 which already held a
# different (non-flag) numeric value and is left untouched.

$wb = $excel.ActiveWorkbook
$phases = $wb.Worksheets.Item("Phases")
$events = $wb.Worksheets.Item("Events")

# --- Data change: Events!I2:I57, flip 1 -> 0 (row 5 is a day-count, not a flag) ---
for ($r = 2; $r -le 57; $r++) {
    $cell = $events.Cells.Item($r, 9)
    if ($cell.Value2 -eq 1) {
        $cell.Value2 = 0
    }
}

# --- View state: the Events tab becomes the active/selected tab, scrolled
#     further down with J62 selected; Phases keeps its existing F6 selection ---
$phases.Range("F6").Select() | Out-Null
$events.Activate() | Out-Null
$events.Range("J62").Select() | Out-Null
